# Update "想去人数" (interested-people count) figures across the
# "展览" (sheet1), "演出" (sheet2) and "全部类型" (sheet4) sheets,
# matching the refreshed data pull recorded in the commit.

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 770
$wsExpo.Range("F6").Value = 125
$wsExpo.Range("F8").Value = 129
$wsExpo.Range("F9").Value = 326
$wsExpo.Range("F10").Value = 439
$wsExpo.Range("F11").Value = 499
$wsExpo.Range("F12").Value = 135
$wsExpo.Range("F13").Value = 11486
$wsExpo.Range("F14").Value = 5380

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F4").Value = 7

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 770
$wsAll.Range("F8").Value = 125
$wsAll.Range("F10").Value = 129
$wsAll.Range("F11").Value = 326
$wsAll.Range("F12").Value = 439
$wsAll.Range("F13").Value = 499
$wsAll.Range("F14").Value = 135
$wsAll.Range("F15").Value = 11486
$wsAll.Range("F16").Value = 7
$wsAll.Range("F17").Value = 5380
